$d = $word.ActiveDocument

# Locate the "Date" paragraph that holds "7/8/2020" so we can insert the
# new paragraph right after it (and before the "Heading2" paragraph).
$dateParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "7/8/2020") {
        $dateParagraph = $p
        break
    }
}

# Insert a brand new paragraph right after the date paragraph.
$dateParagraph.Range.InsertParagraphAfter()

# The freshly inserted paragraph is the one right after the date paragraph.
$newParagraph = $dateParagraph.Next()
$newParagraph.Style = "FirstParagraph"

$text = "{r setup, include=FALSE} knitr::opts_chunk`$set(echo = TRUE)"
$newParagraph.Range.Text = $text

# Apply the "VerbatimChar" character style to the run via Find/Replace,
# which (unlike a direct Range.Style/CharacterStyle assignment) writes a
# proper <w:rStyle> run-property rather than clobbering the paragraph
# style.
$find = $newParagraph.Range.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Replacement.Style = "VerbatimChar"
$find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2)
